# Add a "Save" column (H) to the sheet, matching the formatting of the
# existing header cells (e.g. G1) and add the numeric value for row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save", formatted like the other header cells.
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New data cell H2: numeric value 1.
$ws.Range("H2").Value = 1
